{"js": "// Insert a new \"record\" block right after the paragraph\n// \"Lugar: Frutal, Minas Gerais.\" (and before the separator line that\n// already follows it), matching the authored diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its exact text.\nconst anchorText = \"Lugar: Frutal, Minas Gerais.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// New lines to insert, in document order, right after the anchor.\nconst newLines = [\n  \"--------------------------------------------------------------------\",\n  \"\",\n  \"Modelo: FONTE 120A\",\n  \"URL: https://produto.mercadolivre.com.br/MLB-4830312164-fonte-automotiva-carregador-storm-120-a-jfa-_JM\",\n  \"Nome: Fonte Automotiva Carregador Storm 120-a Jfa\",\n  \"Pre\u00e7o: 599.5\",\n  \"Pre\u00e7o Previsto: 634.4\",\n  \"Loja: MORO1373899\",\n  \"Tipo: Cl\u00e1ssico\",\n  \"Lugar: Rio de Janeiro, Rio de Janeiro.\"\n];\n\n// Insert each line directly after the anchor, one at a time, re-anchoring\n// on the paragraph we just created so the whole block lands in order\n// right after \"Lugar: Frutal, Minas Gerais.\" and before the separator\n// that used to directly follow it.\nlet current = anchor;\nfor (const line of newLines) {\n  current = current.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert a new \"record\" block right after the paragraph\n# \"Lugar: Frutal, Minas Gerais.\" (and before the separator line that\n# already follows it), matching the authored diff.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its exact text.\n$find = $d.Content.Find\n$find.Text = \"Lugar: Frutal, Minas Gerais.\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find anchor paragraph: Lugar: Frutal, Minas Gerais.\"\n}\n$anchorIndex = $find.Parent.Paragraphs(1).Index\n\n# New lines to insert, in document order, right after the anchor.\n$newLines = @(\n    \"--------------------------------------------------------------------\",\n    \"\",\n    \"Modelo: FONTE 120A\",\n    \"URL: https://produto.mercadolivre.com.br/MLB-4830312164-fonte-automotiva-carregador-storm-120-a-jfa-_JM\",\n    \"Nome: Fonte Automotiva Carregador Storm 120-a Jfa\",\n    \"Pre\u00e7o: 599.5\",\n    \"Pre\u00e7o Previsto: 634.4\",\n    \"Loja: MORO1373899\",\n    \"Tipo: Cl\u00e1ssico\",\n    \"Lugar: Rio de Janeiro, Rio de Janeiro.\"\n)\n\n# Insert each line directly after the running index, re-fetching the\n# paragraph object fresh each time (stale Range objects do not reliably\n# track newly-inserted paragraph breaks in this host).\n$insertAfterIdx = $anchorIndex\nforeach ($line in $newLines) {\n    $p = $d.Paragraphs($insertAfterIdx)\n    $p.Range.InsertParagraphAfter()\n    $insertAfterIdx = $insertAfterIdx + 1\n    if ($line -ne \"\") {\n        $d.Paragraphs($insertAfterIdx).Range.Text = $line\n    }\n}\n"}
